$wb = $excel.ActiveWorkbook

# Sheet: 展览 (column F = "想去人数" / "want to go" counts)
$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 9026
$ws.Cells.Item(4, 6).Value = 327
$ws.Cells.Item(6, 6).Value = 751
$ws.Cells.Item(7, 6).Value = 141
$ws.Cells.Item(8, 6).Value = 85
$ws.Cells.Item(9, 6).Value = 390
$ws.Cells.Item(10, 6).Value = 925
$ws.Cells.Item(11, 6).Value = 4081
$ws.Cells.Item(12, 6).Value = 328
$ws.Cells.Item(13, 6).Value = 206
$ws.Cells.Item(14, 6).Value = 822
$ws.Cells.Item(15, 6).Value = 786
$ws.Cells.Item(17, 6).Value = 1
$ws.Cells.Item(18, 6).Value = 514
$ws.Cells.Item(19, 6).Value = 10
$ws.Cells.Item(20, 6).Value = 25
$ws.Cells.Item(21, 6).Value = 1476
$ws.Cells.Item(22, 6).Value = 1378
$ws.Cells.Item(23, 6).Value = 552
$ws.Cells.Item(25, 6).Value = 160
$ws.Cells.Item(26, 6).Value = 190
$ws.Cells.Item(27, 6).Value = 403
$ws.Cells.Item(28, 6).Value = 83
$ws.Cells.Item(29, 6).Value = 1028
$ws.Cells.Item(32, 6).Value = 808
$ws.Cells.Item(33, 6).Value = 87
$ws.Cells.Item(35, 6).Value = 130
$ws.Cells.Item(37, 6).Value = 32
$ws.Cells.Item(40, 6).Value = 216
$ws.Cells.Item(41, 6).Value = 449
$ws.Cells.Item(42, 6).Value = 42
$ws.Cells.Item(43, 6).Value = 36

# Sheet: 演出 (column F = "想去人数" / "want to go" counts)
$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(2, 6).Value = 28
$ws.Cells.Item(4, 6).Value = 108
$ws.Cells.Item(6, 6).Value = 80

# Sheet: 全部类型 (column F = "想去人数" / "want to go" counts)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(3, 6).Value = 9026
$ws.Cells.Item(4, 6).Value = 327
$ws.Cells.Item(5, 6).Value = 751
$ws.Cells.Item(6, 6).Value = 141
$ws.Cells.Item(7, 6).Value = 85
$ws.Cells.Item(8, 6).Value = 390
$ws.Cells.Item(9, 6).Value = 925
$ws.Cells.Item(10, 6).Value = 28
$ws.Cells.Item(11, 6).Value = 4081
$ws.Cells.Item(12, 6).Value = 328
$ws.Cells.Item(13, 6).Value = 206
$ws.Cells.Item(15, 6).Value = 108
$ws.Cells.Item(16, 6).Value = 822
$ws.Cells.Item(17, 6).Value = 786
$ws.Cells.Item(19, 6).Value = 80
$ws.Cells.Item(21, 6).Value = 1
$ws.Cells.Item(22, 6).Value = 514
$ws.Cells.Item(23, 6).Value = 10
$ws.Cells.Item(25, 6).Value = 25
$ws.Cells.Item(26, 6).Value = 1476
$ws.Cells.Item(27, 6).Value = 1378
$ws.Cells.Item(28, 6).Value = 553
$ws.Cells.Item(30, 6).Value = 160
$ws.Cells.Item(31, 6).Value = 190
$ws.Cells.Item(33, 6).Value = 403
$ws.Cells.Item(34, 6).Value = 83
$ws.Cells.Item(35, 6).Value = 1028
$ws.Cells.Item(37, 6).Value = 808
$ws.Cells.Item(38, 6).Value = 87
$ws.Cells.Item(40, 6).Value = 130
$ws.Cells.Item(42, 6).Value = 32
$ws.Cells.Item(44, 6).Value = 216
$ws.Cells.Item(45, 6).Value = 449
$ws.Cells.Item(46, 6).Value = 42
$ws.Cells.Item(47, 6).Value = 36

